# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.018.50'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").Value = '3.060.80'
$ws.Range("E3").Value = '  +2.62%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''578.39'
$ws.Range("E5").Value = '  +0.23%  '

$ws.Range("D6").Value = '''166.85'
$ws.Range("E6").Value = '  +2.40%  '

$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").Value = '3.058.04'
$ws.Range("E8").Value = '  +2.61%  '

$ws.Range("D9").Value = '''0.522'
$ws.Range("E9").Value = '  +1.03%  '

$ws.Range("E10").Value = '  -1.06%  '

$ws.Range("D11").Value = '''0.152'
$ws.Range("E11").Value = '  -0.60%  '

$ws.Range("E12").Value = '  +6.00%  '

$ws.Range("D13").Value = '''0.0000248'
$ws.Range("E13").Value = '  +0.10%  '

$ws.Range("D14").Value = '''36.60'
$ws.Range("E14").Value = '  +5.33%  '

$ws.Range("E15").Value = '  -0.44%  '

$ws.Range("D16").Value = '3.570.74'
$ws.Range("E16").Value = '  +2.71%  '

$ws.Range("D17").Value = '66.047.16'
$ws.Range("E17").Value = '  -0.15%  '

$ws.Range("D18").Value = '''7.17'
$ws.Range("E18").Value = '  +3.59%  '

$ws.Range("D19").Value = '3.061.52'
$ws.Range("E19").Value = '  +2.62%  '

$ws.Range("D20").Value = '''16.21'
$ws.Range("E20").Value = '  +17.08%  '

$ws.Range("D21").Value = '''461.21'
$ws.Range("E21").Value = '  +2.23%  '

$ws.Range("D22").Value = '''0.711'
$ws.Range("E22").Value = '  +4.37%  '

$ws.Range("D23").Value = '''7.40'
$ws.Range("E23").Value = '  +1.54%  '

$ws.Range("D24").Value = '''82.95'
$ws.Range("E24").Value = '  +1.05%  '

$ws.Range("D25").Value = '''12.82'
$ws.Range("E25").Value = '  +4.65%  '

$ws.Range("D26").Value = '''2.26'
$ws.Range("E26").Value = '  +0.69%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '''10.05'
$ws.Range("E27").Value = '  -0.79%  '

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").Value = '''8.09'
$ws.Range("E29").Value = '  -0.87%  '

$ws.Range("D30").Value = '''2.42'
$ws.Range("E30").Value = '  +0.19%  '

$ws.Range("E31").Value = '  +2.03%  '

$ws.Range("D32").Value = '''0.0000102'
$ws.Range("E32").Value = '  +0.14%  '

$ws.Range("D33").Value = '''28.18'
$ws.Range("E33").Value = '  +3.32%  '

$ws.Range("E34").Value = '  +3.79%  '

$ws.Range("D35").Value = '''0.998'
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("D36").Value = '''0.993'
$ws.Range("E36").Value = '  +1.24%  '

$ws.Range("D37").Value = '''5.85'
$ws.Range("E37").Value = '  +1.11%  '

$ws.Range("D38").Value = '''48.26'
$ws.Range("E38").Value = '  +8.63%  '

$ws.Range("D39").Value = '''49.83'
$ws.Range("E39").Value = '  +0.71%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''2.03'
$ws.Range("E40").Value = '  -1.02%  '

$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").Value = '''0.311'
$ws.Range("E41").Value = '  +3.11%  '

$ws.Range("E42").Value = '  +1.65%  '

$ws.Range("D43").Value = '''2.86'
$ws.Range("E43").Value = '  -0.57%  '

$ws.Range("D44").Value = '''8.63'
$ws.Range("E44").Value = '  +2.57%  '

$ws.Range("D45").Value = '''0.0360'
$ws.Range("E45").Value = '  +1.11%  '

$ws.Range("D46").Value = '''379.29'
$ws.Range("E46").Value = '  -3.10%  '

$ws.Range("D47").Value = '2.753.93'
$ws.Range("E47").Value = '  +0.81%  '

$ws.Range("D48").Value = '''133.75'
$ws.Range("E48").Value = '  +1.41%  '

$ws.Range("E49").Value = '  +0.01%  '

$ws.Range("D50").Value = '''24.42'
$ws.Range("E50").Value = '  +4.55%  '

$ws.Range("D51").Value = '''2.22'
$ws.Range("E51").Value = '  +3.45%  '
